$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# This handoff report's GUID-named source/handoff files were regenerated
# (new report run), and the handoff/generate timestamps moved forward.
# Old id:  85e77eb3-7e50-4829-ba4e-c78f10025817
# New id:  496dd6b7-0527-4a7a-9b44-c6f0f1c5e936
# ----------------------------------------------------------------------

$newId = "496dd6b7-0527-4a7a-9b44-c6f0f1c5e936"

$newZhXlf = "$newId.f32e101145f60b56120d6be6fdeea4cd54dbad66.zh-cn.xlf"
$newDeXlf = "$newId.f32e101145f60b56120d6be6fdeea4cd54dbad66.de-de.xlf"

# ---------------- Overview sheet ----------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newId.md"
$wsOverview.Range("B2").Value = "e2e\$newId.md"
$wsOverview.Range("G2").Value = "2016-09-06 17:38:29"

foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = "e2e\$newId.md"
}

# ---------------- zh-cn sheet ----------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "$newId.md"
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = "2016-09-06 17:38:21"

foreach ($hl in $wsZh.Hyperlinks) {
    $hl.TextToDisplay = "$newId.md"
}

# ---------------- de-de sheet ----------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "$newId.md"
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = "2016-09-06 17:38:29"

foreach ($hl in $wsDe.Hyperlinks) {
    $hl.TextToDisplay = "$newId.md"
}
